# "regenerate orders with updates distance/sizes"
#
# The underlying trial-order tables encode distance/size conditions inside
# several text fields (Condition, Filename_Left/Right, Distance, Size).
# This run relabels:
#   D64 -> D69
#   D51 -> D55
#   D80 -> D86
#   S30 -> S31
# wherever they occur as substrings of a cell's text, leaving every other
# value (Face ids, durations, booleans, block/trial indices, ConditionID, …)
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$rows = $used.Rows.Count
$cols = $used.Columns.Count

for ($r = 1; $r -le $rows; $r++) {
  for ($c = 1; $c -le $cols; $c++) {
    $cell = $ws.Cells.Item($r, $c)
    $t = $cell.Text

    if ($t -ne $null -and $t -ne "") {
      if ($t -like "*D64*" -or $t -like "*D51*" -or $t -like "*D80*" -or $t -like "*S30*") {
        $new = $t -replace "D64", "D69"
        $new = $new -replace "D51", "D55"
        $new = $new -replace "D80", "D86"
        $new = $new -replace "S30", "S31"

        if ($new -ne $t) {
          $cell.Value = $new
        }
      }
    }
  }
}
